$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (all cells in this sheet store
# their numbers/text as inline strings, so values are written as text).
$changes = [ordered]@{
    "D2"  = "235.77"
    "D4"  = "5.366"
    "D5"  = "0.05568"
    "D6"  = "3.365"
    "D7"  = "6.459"
    "D8"  = "0.8006"
    "D9"  = "1.040"
    "D10" = "0.1402"
    "D11" = "0.07245"
    "D12" = "0.03189"
    "D13" = "0.02936"
    "D14" = "0.09240"
    "D15" = "0.001667"
    "D16" = "3.250"
    "D17" = "0.04764"

    "B18" = "TigerCash"
    "C18" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D18" = "0.006259"
    "E18" = "17TigerCashTCH"

    "B19" = "HotbitToken"
    "C19" = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
    "D19" = "0.005068"
    "E19" = "18HotbitTokenHTB"

    "B20" = "BitKan"
    "C20" = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
    "D20" = "0.001052"
    "E20" = "19BitKanKAN"

    "B21" = "NitroEx"
    "C21" = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
    "D21" = "0.0001500"
    "E21" = "20NitroExNTX"

    "B22" = "UpBots"
    "C22" = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
    "D22" = "0.0004200"
    "E22" = "21UpBotsUBXT"

    "B23" = "LEO"
    "C23" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D23" = "3.949"
    "E23" = "22LEOLEO"

    "B24" = "BTSEToken"
    "C24" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "D24" = "2.200"
    "E24" = "23BTSETokenBTSE"

    "B25" = "One"
    "C25" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "D25" = "0.01163"
    "E25" = "24OneONEBestin24h"

    "D26" = "0.3289"
    "D27" = "0.1296"

    "D41" = "0.006994"

    "B42" = "CEJI"
    "C42" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "D42" = "0.003500"
    "E42" = "41CEJICEJI"

    "B43" = "BKEXToken"
    "C43" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "D43" = "0.1035"
    "E43" = "42BKEXTokenBKK"

    "D44" = "0.008881"
    "D45" = "0.00005440"
    "D47" = "0.6800"
    "D48" = "0.03268"
    "E48" = "47BOLOBOLOWorstin24h"
    "D49" = "0.00002100"
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $changes[$addr]
}
